$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement table (rows 16-19) lists each worker's document
# number / name alongside the overdue "Periodo Mora" they belong to.
# This update re-sorts the entries so the two records for period 1612
# come first (rows 16-17), followed by the two for period 1701
# (rows 18-19) - i.e. it groups by period while keeping each
# worker's own document number/name pair intact.
#   Row 16: DAMARIS ESTHER BENEDETTY BONFANTE / 1047365859 -> period 1612
#   Row 17: EDITH MARIA ROMERO MARTINEZ      / 1143379924 -> period 1612
#   Row 18: DAMARIS ESTHER BENEDETTY BONFANTE / 1047365859 -> period 1701
#   Row 19: EDITH MARIA ROMERO MARTINEZ      / 1143379924 -> period 1701

$ws.Range("C16").Value = "1047365859"
$ws.Range("D16").Value = "DAMARIS ESTHER BENEDETTY BONFANTE"
$ws.Range("E16").Value = "1612"

$ws.Range("C17").Value = "1143379924"
$ws.Range("D17").Value = "EDITH MARIA ROMERO MARTINEZ"
$ws.Range("E17").Value = "1612"

$ws.Range("C18").Value = "1047365859"
$ws.Range("D18").Value = "DAMARIS ESTHER BENEDETTY BONFANTE"
$ws.Range("E18").Value = "1701"

$ws.Range("C19").Value = "1143379924"
$ws.Range("D19").Value = "EDITH MARIA ROMERO MARTINEZ"
$ws.Range("E19").Value = "1701"
